$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10Nov2019")
$ws.Activate()

# --- Row 20: fill in previously-empty Q/R/S, and add new T/U cells ---
$ws.Range("Q20").Value = 0.61519720109939702
$ws.Range("R20").Value = 0.80117443786838605
$ws.Range("S20").Value = 0.95193928700906905
$ws.Range("T20").Value = 0.048338237030663297
$ws.Range("U20").Value = 0.053086872495026999

# --- Insert a new row at 22 (shifts old row 22.. down by one) ---
$ws.Rows.Item(22).Insert()

$green = 0x50D092

# New row 22 content: "yes/yes/free/add 20 topology" result row (highlighted green)
$ws.Range("A22").Value = "yes"
$ws.Range("B22").Value = "yes"
$ws.Range("C22").Value = "free"
$ws.Range("D22").Value = "add 20 topology"
$ws.Range("U22").Value = 0.045328427185508897

$rngPercent = $ws.Range("F22:H22")
$rngPercent.NumberFormat = "0.00%"

$rngPercent2 = $ws.Range("Q22:S22")
$rngPercent2.Value = $rngPercent2.Value
$rngPercent2.NumberFormat = "0.00%"
$ws.Range("Q22").Value = 0.67126424790124095
$ws.Range("R22").Value = 0.90641066109189905
$ws.Range("S22").Value = 0.97211803796346896

$rngComma = $ws.Range("N22:O22")
$rngComma.NumberFormat = "_(* #,##0.00000_);_(* \(#,##0.00000\);_(* ""-""??_);_(@_)"

$ws.Range("A22:D22").Interior.Color = $green
$ws.Range("U22").Interior.Color = $green
$ws.Range("F22:H22").Interior.Color = $green
$ws.Range("Q22:S22").Interior.Color = $green
$ws.Range("N22:O22").Interior.Color = $green

# Row 23 is the shifted-down former blank spacer row; re-touch its fill so the
# number-format styles stay distinct (matches a reformat of that block).
$ws.Range("F23:H23").Interior.Pattern = -4142
$ws.Range("N23:O23").Interior.Pattern = -4142
$ws.Range("Q23:S23").Interior.Pattern = -4142

# --- View state: scroll / selection ---
$ws.Range("U22").Select()
